$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "63.410.54"
$ws.Range("E2").Value = "  +2.31%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.484.68"
$ws.Range("E3").Value = "  +2.73%  "

# Row 5 - BNB
Set-TextValue "D5" "573.87"
$ws.Range("E5").Value = "  +1.85%  "

# Row 6 - Solana
Set-TextValue "D6" "148.92"
$ws.Range("E6").Value = "  +4.33%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - XRP
Set-TextValue "D8" "0.540"
$ws.Range("E8").Value = "  +1.86%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +4.59%  "

# Row 10 - TRON
$ws.Range("E10").Value = "  +0.47%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +3.92%  "

# Row 12 - Toncoin
Set-TextValue "D12" "5.34"
$ws.Range("E12").Value = "  +2.58%  "

# Row 13 - Avalanche
Set-TextValue "D13" "27.33"
$ws.Range("E13").Value = "  +5.46%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +6.55%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.948.81"
$ws.Range("E15").Value = "  +3.26%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "63.280.60"
$ws.Range("E16").Value = "  +2.24%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.488.15"
$ws.Range("E17").Value = "  +3.05%  "

# Row 18 - Chainlink
Set-TextValue "D18" "11.57"
$ws.Range("E18").Value = "  +2.21%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  +6.50%  "

# Row 20 & 21 - swap Polkadot / BitcoinCash plus updated values
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D20" "328.39"
$ws.Range("E20").Value = "  +1.34%  "

$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D21" "4.23"
$ws.Range("E21").Value = "  +2.48%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.04%  "

# Row 23 - SuiNetwork
$ws.Range("E23").Value = "  +10.68%  "

# Row 24 - Litecoin
Set-TextValue "D24" "67.65"
$ws.Range("E24").Value = "  +1.43%  "

# Row 25 - Bittensor
Set-TextValue "D25" "639.88"
$ws.Range("E25").Value = "  +16.37%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  +13.14%  "

# Row 27 - Aptos
Set-TextValue "D27" "8.77"

# Row 28 - WrappedeETH
$ws.Range("D28").Value = "2.656.50"
$ws.Range("E28").Value = "  +4.66%  "

# Row 29 - Fetch.AI
$ws.Range("E29").Value = "  +9.29%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextValue "D30" "8.46"
$ws.Range("E30").Value = "  +3.66%  "

# Row 31 - Binance-PegBSC-USD
Set-TextValue "D31" "0.999"
$ws.Range("E31").Value = "  -0.14%  "

# Row 32 - Kaspa
$ws.Range("E32").Value = "  -1.90%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  +2.74%  "

# Row 34 - NEARProtocol
$ws.Range("E34").Value = "  +10.31%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  +3.52%  "

# Row 37 - PolygonEcosystemToken
$ws.Range("E37").Value = "  +2.00%  "

# Row 38 - RenderToken
Set-TextValue "D38" "5.51"
$ws.Range("E38").Value = "  +1.76%  "

# Row 39 - EthereumClassic
Set-TextValue "D39" "18.92"
$ws.Range("E39").Value = "  +2.01%  "

# Row 40 - Stacks
Set-TextValue "D40" "1.86"
$ws.Range("E40").Value = "  +3.29%  "

# Row 41 - Monero
Set-TextValue "D41" "147.04"
$ws.Range("E41").Value = "  -4.08%  "

# Row 42 - dogwifhat
Set-TextValue "D42" "2.64"
$ws.Range("E42").Value = "  +18.73%  "

# Row 43 - USDe
$ws.Range("E43").Value = "  +0.80%  "

# Row 44 - Aave
Set-TextValue "D44" "150.51"
$ws.Range("E44").Value = "  +2.63%  "

# Row 45 - Filecoin
$ws.Range("E45").Value = "  +3.65%  "

# Row 46 - InjectiveProtocol
$ws.Range("E46").Value = "  +6.96%  "

# Row 47 - Hedera
$ws.Range("E47").Value = "  +4.28%  "

# Row 48 - Mantle
$ws.Range("E48").Value = "  +2.90%  "

# Row 49 - VeChain
$ws.Range("E49").Value = "  +5.83%  "

# Row 50 - Stellar
Set-TextValue "D50" "0.0929"
$ws.Range("E50").Value = "  +1.02%  "

# Row 51 - ONDO
$ws.Range("E51").Value = "  +5.54%  "
